# Applies the "Saldo_guide" update: bump the reference date from 2024-07-04
# to 2024-07-05, refresh a handful of balance values that changed between
# the two extract runs, rename the sheet to match the new extract
# timestamp, and restore the last-used selection (K12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the worksheet tab to the new extraction timestamp ---------
$ws.Name = "IClientBalance-20240705-091652-"

# --- 2) Bump every "Dt. Referencia" (column G) value by one day ----------
# All data rows (2..275) held the same serial date (45477 = 2024-07-04);
# they all move to 45478 (2024-07-05).
for ($r = 2; $r -le 275; $r++) {
    $ws.Cells.Item($r, 7).Value = 45478
}

# --- 3) Update the handful of rows whose Saldo Previsto / Vl. Total ------
#        values were corrected in this extract ----------------------------
$updatedRows = @{
    109 = 192.47
    110 = 326.64
    158 = 730.49
    161 = 252.11
    230 = 404.28
    231 = 1452.05
}

foreach ($r in $updatedRows.Keys) {
    $val = $updatedRows[$r]
    $ws.Cells.Item($r, 5).Value = $val   # column E - Saldo Previsto
    $ws.Cells.Item($r, 8).Value = $val   # column H - Vl. Total
}

# --- 4) Restore the active cell selection (K12) on the sheet --------------
$ws.Range("K12").Select()
